# Insert 3 new data rows (882-884) on the "Fruta, Terminal La Palmera de La
# Serena - Mandarina" sheet, right above the existing "Murcott" block that
# was previously at rows 882-885 (now pushed down to 885-888 and beyond).
# Net effect: dimension grows from A1:T979 to A1:T982.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 882 onward down by 3 (inserts 3 blank rows, copies formatting
# from the row above - keeps the date style s="2" on column D).
$ws.Range("A882:A884").EntireRow.Insert()

# --- Row 882: Murcott / Especial -------------------------------------------------
$ws.Cells.Item(882,1).Value  = 8
$ws.Cells.Item(882,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(882,3).Value  = "Coquimbo"
$ws.Cells.Item(882,4).Value  = 45194
$ws.Cells.Item(882,5).Value  = 4
$ws.Cells.Item(882,6).Value  = "Fruta"
$ws.Cells.Item(882,7).Value  = 100102
$ws.Cells.Item(882,8).Value  = "Cítricos"
$ws.Cells.Item(882,9).Value  = 100102004
$ws.Cells.Item(882,10).Value = "Mandarina"
$ws.Cells.Item(882,11).Value = "Murcott"
$ws.Cells.Item(882,12).Value = "Especial"
$ws.Cells.Item(882,13).Value = 500
$ws.Cells.Item(882,14).Value = 6500
$ws.Cells.Item(882,15).Value = 7000
$ws.Cells.Item(882,16).Value = 6750
$ws.Cells.Item(882,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(882,18).Value = "Provincia de Limarí"
$ws.Cells.Item(882,19).Value = 675
$ws.Cells.Item(882,20).Value = 10

# --- Row 883: Murcott / Primera --------------------------------------------------
$ws.Cells.Item(883,1).Value  = 8
$ws.Cells.Item(883,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(883,3).Value  = "Coquimbo"
$ws.Cells.Item(883,4).Value  = 45194
$ws.Cells.Item(883,5).Value  = 4
$ws.Cells.Item(883,6).Value  = "Fruta"
$ws.Cells.Item(883,7).Value  = 100102
$ws.Cells.Item(883,8).Value  = "Cítricos"
$ws.Cells.Item(883,9).Value  = 100102004
$ws.Cells.Item(883,10).Value = "Mandarina"
$ws.Cells.Item(883,11).Value = "Murcott"
$ws.Cells.Item(883,12).Value = "Primera"
$ws.Cells.Item(883,13).Value = 600
$ws.Cells.Item(883,14).Value = 5500
$ws.Cells.Item(883,15).Value = 6000
$ws.Cells.Item(883,16).Value = 5750
$ws.Cells.Item(883,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(883,18).Value = "Provincia de Limarí"
$ws.Cells.Item(883,19).Value = 575
$ws.Cells.Item(883,20).Value = 10

# --- Row 884: Murcott / Segunda --------------------------------------------------
$ws.Cells.Item(884,1).Value  = 8
$ws.Cells.Item(884,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(884,3).Value  = "Coquimbo"
$ws.Cells.Item(884,4).Value  = 45194
$ws.Cells.Item(884,5).Value  = 4
$ws.Cells.Item(884,6).Value  = "Fruta"
$ws.Cells.Item(884,7).Value  = 100102
$ws.Cells.Item(884,8).Value  = "Cítricos"
$ws.Cells.Item(884,9).Value  = 100102004
$ws.Cells.Item(884,10).Value = "Mandarina"
$ws.Cells.Item(884,11).Value = "Murcott"
$ws.Cells.Item(884,12).Value = "Segunda"
$ws.Cells.Item(884,13).Value = 600
$ws.Cells.Item(884,14).Value = 4500
$ws.Cells.Item(884,15).Value = 5000
$ws.Cells.Item(884,16).Value = 4750
$ws.Cells.Item(884,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(884,18).Value = "Provincia de Limarí"
$ws.Cells.Item(884,19).Value = 475
$ws.Cells.Item(884,20).Value = 10
